$d = $word.ActiveDocument

# --- Step 1: move the _GoBack bookmark from the end of the document onto the
# "When paying utilities..." paragraph, positioned right after its text (and
# before its paragraph mark), matching the target XML exactly.
#
# The engine has a quirk where Bookmarks.Add() placed on a collapsed range
# whose position equals (paragraph.End - 1) snaps to the wrong location.
# Work around it by temporarily inserting a throwaway character after the
# target position (so the bookmark position is no longer "last char before
# the paragraph mark"), adding the bookmark there, then removing the
# throwaway character again.
$utilPara = $d.Paragraphs.Item(2)
$endPos = $utilPara.Range.End - 1

$tmp = $d.Range($endPos, $endPos)
$tmp.InsertAfter("Z")

$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$d.Range($endPos, $endPos + 1).Delete()

# --- Step 2: remove the "Implement full set" and "Implement trading"
# paragraphs entirely (text + paragraph mark).
$target1 = $d.Content.Find
$target1.ClearFormatting()
$found1 = $d.Content.Find.Execute("Implement full set", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "Implement full set`r") {
        $para.Range.Delete()
        break
    }
}

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "Implement trading`r") {
        $para.Range.Delete()
        break
    }
}
